# fix some item bag bug. optimise the dungeon tips
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dungeon tip text and related data on row 5 (45000002)
$ws.Range("C5").Value = "在12步内击败穷奇"
$ws.Range("D5").Value = 5
$ws.Range("F5").Value = 12

# Row 6 (45000003)
$ws.Range("D6").Value = 3

# Row 7 (45000004)
$ws.Range("D7").Value = 4
$ws.Range("K7").Value = 10

# Rows 8-12 (45000005 .. 45000101) - fill in the "Hard" (difficulty) column D
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("D11").Value = 4
$ws.Range("D12").Value = 5

# Update the active selection shown in the sheet view
$ws.Range("D10").Select()
